$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eventos")

$ws.Range('K2').Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991205821/1748891365835/F7KsVXPM8Y-camera-video-segment-driver-1748891368335.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSFCQDTFFF%2F20250603%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250603T220042Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEEUaCXVzLXdlc3QtMiJIMEYCIQCELgRYOLz6D0EyssQzoCg%2FuixzgCRREtGC31IH5YZehAIhAKEHU9u2sxeKOT2tqfoTLb0k9TmG6GZRU%2FZnrLE1LrHWKt0DCB4QBBoMNzgxMjA0OTQyMjQ0IgyqxcCt6o1aaEhVkdMqugNVmeHRjKSnyAEAIL4zKkYBC3vzXd1o6bOLGzsh%2FIdE25aCqB1gjrqW7mFCt0kiJxBv%2FN6REIhBmUKm8LUgcRAcdbWZDzvobZTiQ6ztDnYuu7Ubo9mQhjz0IBqEEfz9CJhMVwU%2FqC%2Ba1sqW%2FKUNyXDBfIbvzj2Ljk8oLEVKDR8AzO%2FC%2FIQ%2Fco9%2FjYN41irNhk4Fy7BI8Cg0Czf899yL0gPEMAq1mbcCV%2BWNgOYJ0zJAhdBDg%2FeNyTi1uEh%2BiRchTeFDgr%2BhnKoel%2BdMjQ1QZRAO5ypY%2FBMjJCo8ihh0wuCXf3%2B6bDKA%2BjPbiFzsR8UBEeuIjgFDkc%2BxQvLr5I%2BMMX5R%2BVQ6oHI2CZYsJox06%2BzLAroukqqXes76E8eTnzE4CY%2BC3loXy3ccV1aQRUfvJQA63RXnWSbm6oeb2GHocQdDARxTnzF8FLvTuh15bSXqtsbhcvuOhzQ2Ydu2%2Fvo0CBptD7c402nTFNbPAbjW5aaQFnceMfbgIIl6Jp4eXrjtxCMIF8%2FOpZ7YD46oWbv%2BboTpG2NGqphCIsJ11nKYFtK3JjENidJkpMwy8ypcCdDWdwq%2BhWavQNwjnMP5MOW%2B%2FcEGOqQBLgQq6KMW4S4MOGgTr0%2Fe8bbOvGonHumdZ3xyAPDiJX9nHwEHyVOkrle1oHUYoTiXWsaEn4zBZIBW7zsGyrK4le%2FeuaE2ENwAFqw7YObgDLNBQ7YE0c2gVxFP3PENqXZ2mj8udD9zipFcIdrDWOhrhvcIgpQxAIR%2BnoNWOzNujLGatFjuR3yY1LBRnrL9YHBm8rHKE4D4A9AmCb8No%2B1%2BWGiJ3bE%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2004%20Jun%202025%2006%3A00%3A42%20GMT&X-Amz-Signature=021884838202ebd4730d994377b97b8325c6648c4b07f4198671823c73e14ab0'
$ws.Range('L2').Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205821/1748891365835/MSujKlorgb-camera-video-segment-1748891368335.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSFCQDTFFF%2F20250603%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250603T220042Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEEUaCXVzLXdlc3QtMiJIMEYCIQCELgRYOLz6D0EyssQzoCg%2FuixzgCRREtGC31IH5YZehAIhAKEHU9u2sxeKOT2tqfoTLb0k9TmG6GZRU%2FZnrLE1LrHWKt0DCB4QBBoMNzgxMjA0OTQyMjQ0IgyqxcCt6o1aaEhVkdMqugNVmeHRjKSnyAEAIL4zKkYBC3vzXd1o6bOLGzsh%2FIdE25aCqB1gjrqW7mFCt0kiJxBv%2FN6REIhBmUKm8LUgcRAcdbWZDzvobZTiQ6ztDnYuu7Ubo9mQhjz0IBqEEfz9CJhMVwU%2FqC%2Ba1sqW%2FKUNyXDBfIbvzj2Ljk8oLEVKDR8AzO%2FC%2FIQ%2Fco9%2FjYN41irNhk4Fy7BI8Cg0Czf899yL0gPEMAq1mbcCV%2BWNgOYJ0zJAhdBDg%2FeNyTi1uEh%2BiRchTeFDgr%2BhnKoel%2BdMjQ1QZRAO5ypY%2FBMjJCo8ihh0wuCXf3%2B6bDKA%2BjPbiFzsR8UBEeuIjgFDkc%2BxQvLr5I%2BMMX5R%2BVQ6oHI2CZYsJox06%2BzLAroukqqXes76E8eTnzE4CY%2BC3loXy3ccV1aQRUfvJQA63RXnWSbm6oeb2GHocQdDARxTnzF8FLvTuh15bSXqtsbhcvuOhzQ2Ydu2%2Fvo0CBptD7c402nTFNbPAbjW5aaQFnceMfbgIIl6Jp4eXrjtxCMIF8%2FOpZ7YD46oWbv%2BboTpG2NGqphCIsJ11nKYFtK3JjENidJkpMwy8ypcCdDWdwq%2BhWavQNwjnMP5MOW%2B%2FcEGOqQBLgQq6KMW4S4MOGgTr0%2Fe8bbOvGonHumdZ3xyAPDiJX9nHwEHyVOkrle1oHUYoTiXWsaEn4zBZIBW7zsGyrK4le%2FeuaE2ENwAFqw7YObgDLNBQ7YE0c2gVxFP3PENqXZ2mj8udD9zipFcIdrDWOhrhvcIgpQxAIR%2BnoNWOzNujLGatFjuR3yY1LBRnrL9YHBm8rHKE4D4A9AmCb8No%2B1%2BWGiJ3bE%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2004%20Jun%202025%2006%3A00%3A42%20GMT&X-Amz-Signature=9b23d3781d844e5e6cf38422f7abae45252719be519697e9c00f1c7d39b73e71'
$ws.Range('K3').Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991205821/1748885144009/VaUCDIqpi8-camera-video-segment-driver-1748885146509.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSFCQDTFFF%2F20250603%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250603T220042Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEEUaCXVzLXdlc3QtMiJIMEYCIQCELgRYOLz6D0EyssQzoCg%2FuixzgCRREtGC31IH5YZehAIhAKEHU9u2sxeKOT2tqfoTLb0k9TmG6GZRU%2FZnrLE1LrHWKt0DCB4QBBoMNzgxMjA0OTQyMjQ0IgyqxcCt6o1aaEhVkdMqugNVmeHRjKSnyAEAIL4zKkYBC3vzXd1o6bOLGzsh%2FIdE25aCqB1gjrqW7mFCt0kiJxBv%2FN6REIhBmUKm8LUgcRAcdbWZDzvobZTiQ6ztDnYuu7Ubo9mQhjz0IBqEEfz9CJhMVwU%2FqC%2Ba1sqW%2FKUNyXDBfIbvzj2Ljk8oLEVKDR8AzO%2FC%2FIQ%2Fco9%2FjYN41irNhk4Fy7BI8Cg0Czf899yL0gPEMAq1mbcCV%2BWNgOYJ0zJAhdBDg%2FeNyTi1uEh%2BiRchTeFDgr%2BhnKoel%2BdMjQ1QZRAO5ypY%2FBMjJCo8ihh0wuCXf3%2B6bDKA%2BjPbiFzsR8UBEeuIjgFDkc%2BxQvLr5I%2BMMX5R%2BVQ6oHI2CZYsJox06%2BzLAroukqqXes76E8eTnzE4CY%2BC3loXy3ccV1aQRUfvJQA63RXnWSbm6oeb2GHocQdDARxTnzF8FLvTuh15bSXqtsbhcvuOhzQ2Ydu2%2Fvo0CBptD7c402nTFNbPAbjW5aaQFnceMfbgIIl6Jp4eXrjtxCMIF8%2FOpZ7YD46oWbv%2BboTpG2NGqphCIsJ11nKYFtK3JjENidJkpMwy8ypcCdDWdwq%2BhWavQNwjnMP5MOW%2B%2FcEGOqQBLgQq6KMW4S4MOGgTr0%2Fe8bbOvGonHumdZ3xyAPDiJX9nHwEHyVOkrle1oHUYoTiXWsaEn4zBZIBW7zsGyrK4le%2FeuaE2ENwAFqw7YObgDLNBQ7YE0c2gVxFP3PENqXZ2mj8udD9zipFcIdrDWOhrhvcIgpQxAIR%2BnoNWOzNujLGatFjuR3yY1LBRnrL9YHBm8rHKE4D4A9AmCb8No%2B1%2BWGiJ3bE%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2004%20Jun%202025%2006%3A00%3A42%20GMT&X-Amz-Signature=9c64d18788f0915e0831285bea8e6b1a5c2c84e94c5079463ed749f96a5f6c90'
$ws.Range('L3').Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205821/1748885144009/8zw63V5M8l-camera-video-segment-1748885146509.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSFCQDTFFF%2F20250603%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250603T220042Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEEUaCXVzLXdlc3QtMiJIMEYCIQCELgRYOLz6D0EyssQzoCg%2FuixzgCRREtGC31IH5YZehAIhAKEHU9u2sxeKOT2tqfoTLb0k9TmG6GZRU%2FZnrLE1LrHWKt0DCB4QBBoMNzgxMjA0OTQyMjQ0IgyqxcCt6o1aaEhVkdMqugNVmeHRjKSnyAEAIL4zKkYBC3vzXd1o6bOLGzsh%2FIdE25aCqB1gjrqW7mFCt0kiJxBv%2FN6REIhBmUKm8LUgcRAcdbWZDzvobZTiQ6ztDnYuu7Ubo9mQhjz0IBqEEfz9CJhMVwU%2FqC%2Ba1sqW%2FKUNyXDBfIbvzj2Ljk8oLEVKDR8AzO%2FC%2FIQ%2Fco9%2FjYN41irNhk4Fy7BI8Cg0Czf899yL0gPEMAq1mbcCV%2BWNgOYJ0zJAhdBDg%2FeNyTi1uEh%2BiRchTeFDgr%2BhnKoel%2BdMjQ1QZRAO5ypY%2FBMjJCo8ihh0wuCXf3%2B6bDKA%2BjPbiFzsR8UBEeuIjgFDkc%2BxQvLr5I%2BMMX5R%2BVQ6oHI2CZYsJox06%2BzLAroukqqXes76E8eTnzE4CY%2BC3loXy3ccV1aQRUfvJQA63RXnWSbm6oeb2GHocQdDARxTnzF8FLvTuh15bSXqtsbhcvuOhzQ2Ydu2%2Fvo0CBptD7c402nTFNbPAbjW5aaQFnceMfbgIIl6Jp4eXrjtxCMIF8%2FOpZ7YD46oWbv%2BboTpG2NGqphCIsJ11nKYFtK3JjENidJkpMwy8ypcCdDWdwq%2BhWavQNwjnMP5MOW%2B%2FcEGOqQBLgQq6KMW4S4MOGgTr0%2Fe8bbOvGonHumdZ3xyAPDiJX9nHwEHyVOkrle1oHUYoTiXWsaEn4zBZIBW7zsGyrK4le%2FeuaE2ENwAFqw7YObgDLNBQ7YE0c2gVxFP3PENqXZ2mj8udD9zipFcIdrDWOhrhvcIgpQxAIR%2BnoNWOzNujLGatFjuR3yY1LBRnrL9YHBm8rHKE4D4A9AmCb8No%2B1%2BWGiJ3bE%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2004%20Jun%202025%2006%3A00%3A42%20GMT&X-Amz-Signature=60f868d3f22412ab647cbfc456f4852b7215026c38d12cb8dd1bddff30bfd9ca'
